# Updated cryptos list (price + 1h volume change columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.097.21"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "2.604.88"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").Value = "150.05"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "2.603.76"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  +2.74%  "
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").Value = "27.17"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").Value = "3.077.40"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").Value = "66.952.59"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "2.601.07"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "367.73"
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("D20").Value = "11.03"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").Value = "7.38"
$ws.Range("E21").Value = "  -3.28%  "
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("D25").Value = "73.34"
$ws.Range("E25").Value = "  +4.65%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "9.92"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D29").Value = "584.09"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("D31").Value = "0.0₃0985"
$ws.Range("E31").Value = "  -6.64%  "
$ws.Range("E32").Value = "  -5.21%  "
$ws.Range("D33").Value = "7.64"
$ws.Range("E33").Value = "  -3.60%  "
$ws.Range("E34").Value = "  -3.24%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("E37").Value = "  -2.86%  "
$ws.Range("D38").Value = "156.39"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "19.03"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("D40").Value = "0.364"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  -3.33%  "
$ws.Range("E43").Value = "  -4.34%  "
$ws.Range("D44").Value = "'17.10"
$ws.Range("E44").Value = "  +4.18%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "152.76"
$ws.Range("E46").Value = "  -2.42%  "
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("E49").Value = "  -3.63%  "
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("D51").Value = "21.39"
$ws.Range("E51").Value = "  +1.79%  "
